{"js": "// Replace the date line and all the two-digit multiplication problems\n// with the values from the new day's worksheet. Every \"old\" string in\n// this document is unique, so a direct search + replace per pair is\n// unambiguous and independent of ordering.\nconst replacements = [\n  [\"2026-01-12 Monday\", \"2026-01-13 Tuesday\"],\n  [\"99\u00d752=\", \"53\u00d759=\"],\n  [\"18\u00d746=\", \"99\u00d799=\"],\n  [\"75\u00d719=\", \"43\u00d713=\"],\n  [\"72\u00d735=\", \"31\u00d777=\"],\n  [\"40\u00d768=\", \"72\u00d727=\"],\n  [\"14\u00d713=\", \"22\u00d720=\"],\n  [\"81\u00d790=\", \"16\u00d798=\"],\n  [\"78\u00d734=\", \"48\u00d787=\"],\n  [\"55\u00d712=\", \"76\u00d782=\"],\n  [\"80\u00d719=\", \"95\u00d782=\"],\n  [\"50\u00d746=\", \"68\u00d780=\"],\n  [\"19\u00d762=\", \"26\u00d741=\"],\n  [\"14\u00d772=\", \"25\u00d784=\"],\n  [\"13\u00d732=\", \"30\u00d741=\"],\n  [\"50\u00d792=\", \"65\u00d751=\"],\n  [\"97\u00d714=\", \"44\u00d725=\"],\n  [\"79\u00d734=\", \"47\u00d733=\"],\n  [\"19\u00d720=\", \"15\u00d714=\"],\n  [\"97\u00d777=\", \"67\u00d772=\"],\n  [\"96\u00d784=\", \"87\u00d747=\"],\n  [\"62\u00d729=\", \"17\u00d778=\"],\n  [\"45\u00d786=\", \"54\u00d797=\"],\n  [\"85\u00d713=\", \"35\u00d774=\"],\n  [\"56\u00d765=\", \"24\u00d752=\"],\n  [\"79\u00d760=\", \"20\u00d779=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and all the two-digit multiplication problems\n# with the values from the new day's worksheet. Every \"old\" string in\n# this document is unique, so a direct Find/Replace per pair is\n# unambiguous and independent of ordering.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2026-01-12 Monday\", \"2026-01-13 Tuesday\"),\n    @(\"99\u00d752=\", \"53\u00d759=\"),\n    @(\"18\u00d746=\", \"99\u00d799=\"),\n    @(\"75\u00d719=\", \"43\u00d713=\"),\n    @(\"72\u00d735=\", \"31\u00d777=\"),\n    @(\"40\u00d768=\", \"72\u00d727=\"),\n    @(\"14\u00d713=\", \"22\u00d720=\"),\n    @(\"81\u00d790=\", \"16\u00d798=\"),\n    @(\"78\u00d734=\", \"48\u00d787=\"),\n    @(\"55\u00d712=\", \"76\u00d782=\"),\n    @(\"80\u00d719=\", \"95\u00d782=\"),\n    @(\"50\u00d746=\", \"68\u00d780=\"),\n    @(\"19\u00d762=\", \"26\u00d741=\"),\n    @(\"14\u00d772=\", \"25\u00d784=\"),\n    @(\"13\u00d732=\", \"30\u00d741=\"),\n    @(\"50\u00d792=\", \"65\u00d751=\"),\n    @(\"97\u00d714=\", \"44\u00d725=\"),\n    @(\"79\u00d734=\", \"47\u00d733=\"),\n    @(\"19\u00d720=\", \"15\u00d714=\"),\n    @(\"97\u00d777=\", \"67\u00d772=\"),\n    @(\"96\u00d784=\", \"87\u00d747=\"),\n    @(\"62\u00d729=\", \"17\u00d778=\"),\n    @(\"45\u00d786=\", \"54\u00d797=\"),\n    @(\"85\u00d713=\", \"35\u00d774=\"),\n    @(\"56\u00d765=\", \"24\u00d752=\"),\n    @(\"79\u00d760=\", \"20\u00d779=\")\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
